$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: duplicate row 248 in place (Copy + Insert) -----------------
# This pushes the existing rows 248..261 down to 249..262, and the new
# row 248 starts out as an exact copy of the original row 248.
$ws.Rows.Item(248).Copy()
$ws.Rows.Item(248).Insert()

# --- Step 2: duplicate row 261 (which, after the step-1 shift, holds the
# data that used to be row 260) and insert it at row 262. This pushes the
# row currently at 262 (originally row 261) down to row 263, and leaves a
# fresh copy of the "old row 260" data at row 262. -----------------------
$ws.Rows.Item(261).Copy()
$ws.Rows.Item(262).Insert()

# --- Step 3: rewrite the three fields that differ on the (new) row 248 --
$ws.Cells.Item(248, 4).Value = 44931
$ws.Cells.Item(248, 13).Value = 290
$ws.Cells.Item(248, 18).Value = "Región de O'Higgins"
